$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row text stays logically the same but A1's underlying shared
# string changes from the placeholder "normal_flow_rate_details1" to
# "Operation phase" (B1/C1 keep their existing text).
$ws.Range("A1").Value = "Operation phase"
$ws.Range("B1").Value = "Mean flow rate"
$ws.Range("C1").Value = "Unit"

# Update the "Mean flow rate" column (B2:B9) with the new values.
$ws.Range("B2").Value = 21.509545652562473
$ws.Range("B3").Value = 5.4887138399957749
$ws.Range("B4").Value = 16.180644693854031
$ws.Range("B5").Value = 22.37187719861836
$ws.Range("B6").Value = 4.9961956554845965
$ws.Range("B7").Value = 27.436120310086068
$ws.Range("B8").Value = 30.451713659305238
$ws.Range("B9").Value = 15.822140733987464
